# Weekly refresh of "Fruta, Femacal de La Calera - Caqui" price records.
# The underlying observations (Fecha/Calidad/Volumen/Precio.../Origen) were
# reshuffled across rows 2-32 by the source extract; this applies the new
# per-row values cell by cell, matching the published diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44301
$ws.Range("M2").Value = 45
$ws.Range("D3").Value = 44329
$ws.Range("M3").Value = 56
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 9000
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 900
$ws.Range("D4").Value = 44329
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 8000
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 8000
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 800
$ws.Range("D5").Value = 44322
$ws.Range("M5").Value = 56
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("S5").Value = 1000
$ws.Range("D6").Value = 44322
$ws.Range("M6").Value = 40
$ws.Range("D7").Value = 44319
$ws.Range("M7").Value = 68
$ws.Range("D8").Value = 44319
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 57
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 800
$ws.Range("D9").Value = 44312
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("S9").Value = 1000
$ws.Range("D10").Value = 44315
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 1000
$ws.Range("D11").Value = 44333
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 58
$ws.Range("D12").Value = 44333
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 65
$ws.Range("N12").Value = 9000
$ws.Range("O12").Value = 9000
$ws.Range("P12").Value = 9000
$ws.Range("S12").Value = 900
$ws.Range("D13").Value = 44333
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("S13").Value = 800
$ws.Range("D14").Value = 44326
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 65
$ws.Range("N14").Value = 10000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 10000
$ws.Range("S14").Value = 1000
$ws.Range("D15").Value = 44326
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 67
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("R15").Value = "Provincia de Quillota"
$ws.Range("S15").Value = 800
$ws.Range("D16").Value = 44323
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("R16").Value = "Provincia de Quillota"
$ws.Range("S16").Value = 1000
$ws.Range("D17").Value = 44323
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 9000
$ws.Range("O17").Value = 9000
$ws.Range("P17").Value = 9000
$ws.Range("R17").Value = "Provincia de Quillota"
$ws.Range("S17").Value = 900
$ws.Range("D18").Value = 44307
$ws.Range("M18").Value = 40
$ws.Range("D19").Value = 44321
$ws.Range("M19").Value = 58
$ws.Range("D20").Value = 44302
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("S21").Value = 1000
$ws.Range("D22").Value = 44308
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 48
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("S22").Value = 800
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 56
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 12000
$ws.Range("S23").Value = 1200
$ws.Range("D24").Value = 44699
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 10000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 10000
$ws.Range("R24").Value = "Provincia de Quillota"
$ws.Range("S24").Value = 1000
$ws.Range("D25").Value = 44309
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 45
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("R25").Value = "Provincia de Quillota"
$ws.Range("S25").Value = 1000
$ws.Range("D26").Value = 44314
$ws.Range("M26").Value = 47
$ws.Range("N26").Value = 9000
$ws.Range("O26").Value = 9000
$ws.Range("P26").Value = 9000
$ws.Range("S26").Value = 900
$ws.Range("D27").Value = 44343
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 47
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 1000
$ws.Range("D28").Value = 44343
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 9000
$ws.Range("P28").Value = 9000
$ws.Range("R28").Value = "Región Metropolitana"
$ws.Range("S28").Value = 900
$ws.Range("D29").Value = 44343
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 58
$ws.Range("N29").Value = 8000
$ws.Range("O29").Value = 8000
$ws.Range("P29").Value = 8000
$ws.Range("R29").Value = "Región Metropolitana"
$ws.Range("S29").Value = 800
$ws.Range("D30").Value = 44328
$ws.Range("M30").Value = 45
$ws.Range("N30").Value = 8000
$ws.Range("O30").Value = 8000
$ws.Range("P30").Value = 8000
$ws.Range("S30").Value = 800
$ws.Range("D31").Value = 44328
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 48
$ws.Range("N31").Value = 7000
$ws.Range("O31").Value = 7000
$ws.Range("P31").Value = 7000
$ws.Range("S31").Value = 700
$ws.Range("D32").Value = 44306
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 45
$ws.Range("N32").Value = 10000
$ws.Range("O32").Value = 10000
$ws.Range("P32").Value = 10000
$ws.Range("S32").Value = 1000
